# Update "想去人数" (want-to-go count) figures in the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 12620
$wsExhibit.Range("F3").Value = 603
$wsExhibit.Range("F5").Value = 14
$wsExhibit.Range("F6").Value = 272
$wsExhibit.Range("F7").Value = 391
$wsExhibit.Range("F9").Value = 12593
$wsExhibit.Range("F10").Value = 15
$wsExhibit.Range("F11").Value = 3116
$wsExhibit.Range("F17").Value = 25
$wsExhibit.Range("F21").Value = 6101
$wsExhibit.Range("F23").Value = 3608

$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 12620
$wsAll.Range("F3").Value = 603
$wsAll.Range("F5").Value = 14
$wsAll.Range("F6").Value = 272
$wsAll.Range("F8").Value = 391
$wsAll.Range("F10").Value = 12593
$wsAll.Range("F11").Value = 15
$wsAll.Range("F12").Value = 3116
$wsAll.Range("F18").Value = 25
$wsAll.Range("F23").Value = 6101
$wsAll.Range("F25").Value = 3608
